$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.211
$ws.Range("C2").Value = 8.964
$ws.Range("D2").Value = 23.766
$ws.Range("E2").Value = 8.71
$ws.Range("F2").Value = 7.802
$ws.Range("G2").Value = 18.363
$ws.Range("H2").Value = -5.885
$ws.Range("I2").Value = 0.792
$ws.Range("J2").Value = 12.151
$ws.Range("K2").Value = 17.342
$ws.Range("L2").Value = 4.389
$ws.Range("B3").Value = -14.854
$ws.Range("C3").Value = -10.583
$ws.Range("D3").Value = -22.484
$ws.Range("E3").Value = -15.176
$ws.Range("F3").Value = -20.35
$ws.Range("G3").Value = -9.969
$ws.Range("H3").Value = -41.747
$ws.Range("I3").Value = -25.836
$ws.Range("J3").Value = -2.001
$ws.Range("K3").Value = -6.967
$ws.Range("L3").Value = -16.002
$ws.Range("B4").Value = -13.865
$ws.Range("C4").Value = -11.698
$ws.Range("D4").Value = -27.619
$ws.Range("E4").Value = -11.852
$ws.Range("F4").Value = -18.457
$ws.Range("G4").Value = -8.442
$ws.Range("H4").Value = -39.783
$ws.Range("I4").Value = -21.646
$ws.Range("J4").Value = -3.654
$ws.Range("K4").Value = -3.397
$ws.Range("L4").Value = -14.992
$ws.Range("B5").Value = -16.776
$ws.Range("C5").Value = -6.244
$ws.Range("D5").Value = -33.389
$ws.Range("E5").Value = 3.789
$ws.Range("F5").Value = -14.057
$ws.Range("G5").Value = 3.085
$ws.Range("H5").Value = -32.447
$ws.Range("I5").Value = -9.157
$ws.Range("J5").Value = 3.088
$ws.Range("K5").Value = 4.097
$ws.Range("L5").Value = -21.774
$ws.Range("B6").Value = -11.972
$ws.Range("C6").Value = 41.552
$ws.Range("D6").Value = -33.852
$ws.Range("E6").Value = 47.72
$ws.Range("F6").Value = 32.347
$ws.Range("G6").Value = 18.699
$ws.Range("H6").Value = -36.761
$ws.Range("I6").Value = 0.484
$ws.Range("J6").Value = 18.361
$ws.Range("K6").Value = 76.773
$ws.Range("L6").Value = -20.755
$ws.Range("B7").Value = -2.754
$ws.Range("C7").Value = 7.906
$ws.Range("D7").Value = -8.651
$ws.Range("E7").Value = 8.918
$ws.Range("F7").Value = 6.329
$ws.Range("G7").Value = 3.825
$ws.Range("H7").Value = -9.547
$ws.Range("I7").Value = 0.106
$ws.Range("J7").Value = 3.76
$ws.Range("K7").Value = 13.285
$ws.Range("L7").Value = -4.966
$ws.Range("B8").Value = 93452.13
$ws.Range("C8").Value = 175601.61
$ws.Range("D8").Value = 65999.86
$ws.Range("E8").Value = 69123.93
$ws.Range("F8").Value = 154530.06
$ws.Range("G8").Value = 135809.48
$ws.Range("H8").Value = 74930.63
$ws.Range("I8").Value = 58961
$ws.Range("J8").Value = 43162.4
$ws.Range("K8").Value = 218295.21
$ws.Range("L8").Value = 76172.96
$ws.Range("B9").Value = 0.904
$ws.Range("C9").Value = 0.984
$ws.Range("D9").Value = 1.042
$ws.Range("E9").Value = 0.975
$ws.Range("F9").Value = 1.006
$ws.Range("G9").Value = 1.028
$ws.Range("H9").Value = 1.117
$ws.Range("I9").Value = 0.86
$ws.Range("J9").Value = 0.609
$ws.Range("K9").Value = 1.035
$ws.Range("L9").Value = 0.61
$ws.Range("B10").Value = -10.044
$ws.Range("C10").Value = 0.021
$ws.Range("D10").Value = -16.965
$ws.Range("E10").Value = 1.104
$ws.Range("F10").Value = -1.718
$ws.Range("G10").Value = -4.386
$ws.Range("H10").Value = -18.413
$ws.Range("I10").Value = -6.857
$ws.Range("J10").Value = -1.344
$ws.Range("K10").Value = 5.028
$ws.Range("L10").Value = -10.081
$ws.Range("B11").Value = -0.165
$ws.Range("C11").Value = 0.362
$ws.Range("D11").Value = -0.287
$ws.Range("E11").Value = 0.342
$ws.Range("F11").Value = 0.261
$ws.Range("G11").Value = 0.133
$ws.Range("H11").Value = -0.371
$ws.Range("I11").Value = -0.02
$ws.Range("J11").Value = 0.223
$ws.Range("K11").Value = 0.534
$ws.Range("L11").Value = -0.252
$ws.Range("B12").Value = -0.037
$ws.Range("C12").Value = 0.074
$ws.Range("D12").Value = -0.089
$ws.Range("E12").Value = 0.085
$ws.Range("F12").Value = 0.057
$ws.Range("G12").Value = 0.031
$ws.Range("H12").Value = -0.091
$ws.Range("I12").Value = -0.006
$ws.Range("J12").Value = 0.052
$ws.Range("K12").Value = 0.123
$ws.Range("L12").Value = -0.091
$ws.Range("B13").Value = -49.005
$ws.Range("C13").Value = -32.427
$ws.Range("D13").Value = -66.365
$ws.Range("E13").Value = -44.101
$ws.Range("F13").Value = -39.541
$ws.Range("G13").Value = -36.075
$ws.Range("H13").Value = -50.385
$ws.Range("I13").Value = -42.377
$ws.Range("J13").Value = -22.737
$ws.Range("K13").Value = -34.126
$ws.Range("L13").Value = -34.779
$ws.Range("B14").Value = 20.313
$ws.Range("C14").Value = 20.176
$ws.Range("D14").Value = 32.285
$ws.Range("E14").Value = 24.351
$ws.Range("F14").Value = 21.971
$ws.Range("G14").Value = 24.283
$ws.Range("H14").Value = 27.347
$ws.Range("I14").Value = 24.68
$ws.Range("J14").Value = 14.195
$ws.Range("K14").Value = 23.764
$ws.Range("L14").Value = 22.098
$ws.Range("B15").Value = 0.422
$ws.Range("C15").Value = 0.754
$ws.Range("D15").Value = 0.462
$ws.Range("E15").Value = 0.406
$ws.Range("F15").Value = 0.784
$ws.Range("G15").Value = 0.513
$ws.Range("H15").Value = 0.466
$ws.Range("I15").Value = 0.315
$ws.Range("J15").Value = 0.577
$ws.Range("K15").Value = 0.742
$ws.Range("L15").Value = 0.877
$ws.Range("B16").Value = 7.29
$ws.Range("C16").Value = 7.884
$ws.Range("D16").Value = 8.313
$ws.Range("E16").Value = 7.814
$ws.Range("F16").Value = 8.047
$ws.Range("G16").Value = 8.21
$ws.Range("H16").Value = 8.867
$ws.Range("I16").Value = 6.963
$ws.Range("J16").Value = 5.104
$ws.Range("K16").Value = 8.257
$ws.Range("L16").Value = 5.115